$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 16 (old row 16 "Lag-Llama" and everything
# below shifts down by 4, formulas auto-adjust).
$ws.Rows("16:19").Insert()

# New row 16: UniTime_8 / zero-shot / univariate summary (value only)
$ws.Range("A16").Value = "UniTime_8"
$ws.Range("B16").Value = "zero-shot"
$ws.Range("C16").Value = "univariate"
$ws.Range("M16").Value = 2254

# New row 17: UniTime_8 / finetune / univariate
$ws.Range("A17").Value = "UniTime_8"
$ws.Range("B17").Value = "finetune"
$ws.Range("C17").Value = "univariate"
$ws.Range("D17").Value = 374
$ws.Range("E17").Value = 368
$ws.Range("F17").Value = 1325
$ws.Range("G17").Value = 1320
$ws.Range("H17").Value = 211
$ws.Range("I17").Value = 1282
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 67
$ws.Range("L17").Value = 473
$ws.Range("M17").Formula = "=SUM(D17:L17)"

# New row 18: Unitime_8 / zero-shot / multivariate summary (value only)
$ws.Range("A18").Value = "Unitime_8"
$ws.Range("B18").Value = "zero-shot"
$ws.Range("C18").Value = "multivariate"
$ws.Range("M18").Value = 2877

# New row 19: UniTime_8 / finetune / multivariate
$ws.Range("A19").Value = "UniTime_8"
$ws.Range("B19").Value = "finetune"
$ws.Range("C19").Value = "multivariate"
$ws.Range("D19").Value = 1057
$ws.Range("E19").Value = 1050
$ws.Range("F19").Value = 4320
$ws.Range("G19").Value = 4387
$ws.Range("H19").Value = 632
$ws.Range("I19").Value = 8172
$ws.Range("J19").Value = 99
$ws.Range("K19").Value = 140
$ws.Range("L19").Value = 1564
$ws.Range("M19").Formula = "=SUM(D19:L19)"

# Update view: scroll so row 7 is the top-left visible row, and select O18
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("O18").Select()
